# Update cryptocurrency price/volume data per latest Coinranking pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '79.637.30'
$ws.Range("E2").Value = '  +4.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.146.49'
$ws.Range("E3").Value = '  +2.22%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.19'
$ws.Range("E5").Value = '  +3.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '622.10'
$ws.Range("E6").Value = '  +0.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.264'
$ws.Range("E7").Value = '  +23.34%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.581'
$ws.Range("E9").Value = '  +5.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.147.59'
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.580'
$ws.Range("E11").Value = '  +29.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000248'
$ws.Range("E12").Value = '  +26.34%  '
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.710.47'
$ws.Range("E14").Value = '  +1.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.20'
$ws.Range("E15").Value = '  -1.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.30'
$ws.Range("E16").Value = '  +6.57%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.409.43'
$ws.Range("E17").Value = '  +4.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.133.46'
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.10'
$ws.Range("E19").Value = '  +4.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.94'
$ws.Range("E20").Value = '  +13.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '431.81'
$ws.Range("E21").Value = '  +12.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.01'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.15'
$ws.Range("E23").Value = '  +13.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.90'
$ws.Range("E24").Value = '  +7.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.298.09'
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '75.22'
$ws.Range("E26").Value = '  +3.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.63'
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.71'
$ws.Range("E28").Value = '  +6.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000120'
$ws.Range("E30").Value = '  +10.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.86'
$ws.Range("E32").Value = '  +6.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '550.28'
$ws.Range("E33").Value = '  +9.69%  '
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.148'
$ws.Range("E35").Value = '  +16.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.97'
$ws.Range("E36").Value = '  +2.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.69'
$ws.Range("E37").Value = '  +8.45%  '
$ws.Range("E38").Value = '  +17.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.399'
$ws.Range("E40").Value = '  +5.30%  '
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '20.74'
$ws.Range("E41").Value = '  +3.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '162.84'
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.51'
$ws.Range("E44").Value = '  +6.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '187.78'
$ws.Range("E45").Value = '  -3.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.79'
$ws.Range("E46").Value = '  +7.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.67'
$ws.Range("E47").Value = '  +8.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.778'
$ws.Range("E48").Value = '  -3.43%  '
$ws.Range("E49").Value = '  +1.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.52'
$ws.Range("E50").Value = '  +4.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.19'
$ws.Range("E51").Value = '  +6.54%  '
